$wb = $excel.ActiveWorkbook

# --- Sheet 1: Dades_Període ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("E2").NumberFormat = "@"
$ws1.Range("E2").Value = "09:30 - 10:00"
$ws1.Range("E2").Style = "Normal"
$ws1.Range("H2").NumberFormat = "@"
$ws1.Range("H2").Value = "2026-02-20 10:10:55"
$ws1.Range("H2").Style = "Normal"
$ws1.Range("I2").NumberFormat = "@"
$ws1.Range("I2").Value = "09:30"
$ws1.Range("I2").Style = "Normal"
$ws1.Range("J2").NumberFormat = "@"
$ws1.Range("J2").Value = "https://www.meteo.cat/observacions/xema/dades?codi=XJ&dia=2026-02-20T09:30Z"
$ws1.Range("J2").Style = "Normal"
$ws1.Range("M2").NumberFormat = "@"
$ws1.Range("M2").Value = "186"
$ws1.Range("M2").Style = "Normal"
$ws1.Range("N2").NumberFormat = "@"
$ws1.Range("N2").Value = "50"
$ws1.Range("N2").Style = "Normal"
$ws1.Range("O2").NumberFormat = "@"
$ws1.Range("O2").Value = "1024.2"
$ws1.Range("O2").Style = "Normal"
$ws1.Range("P2").NumberFormat = "@"
$ws1.Range("P2").Value = "0.0"
$ws1.Range("P2").Style = "Normal"
$ws1.Range("Q2").NumberFormat = "@"
$ws1.Range("Q2").Value = "09:30 - 10:00"
$ws1.Range("Q2").Style = "Normal"
$ws1.Range("R2").NumberFormat = "@"
$ws1.Range("R2").Value = "437"
$ws1.Range("R2").Style = "Normal"
$ws1.Range("S2").NumberFormat = "@"
$ws1.Range("S2").Value = "12.8"
$ws1.Range("S2").Style = "Normal"
$ws1.Range("T2").NumberFormat = "@"
$ws1.Range("T2").Value = "11.9"
$ws1.Range("T2").Style = "Normal"
$ws1.Range("U2").NumberFormat = "@"
$ws1.Range("U2").Value = "13.5"
$ws1.Range("U2").Style = "Normal"
$ws1.Range("V2").NumberFormat = "@"
$ws1.Range("V2").Value = "1.4"
$ws1.Range("V2").Style = "Normal"
$ws1.Range("W2").NumberFormat = "@"
$ws1.Range("W2").Value = "5.4"
$ws1.Range("W2").Style = "Normal"
$ws1.Range("X2").NumberFormat = "@"
$ws1.Range("X2").Value = "09:30 - 10:00"
$ws1.Range("X2").Style = "Normal"
$ws1.Range("Y2").NumberFormat = "@"
$ws1.Range("Y2").Value = "12.8"
$ws1.Range("Y2").Style = "Normal"
$ws1.Range("Z2").NumberFormat = "@"
$ws1.Range("Z2").Value = "13.5"
$ws1.Range("Z2").Style = "Normal"
$ws1.Range("AA2").NumberFormat = "@"
$ws1.Range("AA2").Value = "11.9"
$ws1.Range("AA2").Style = "Normal"
$ws1.Range("AB2").NumberFormat = "@"
$ws1.Range("AB2").Value = "50"
$ws1.Range("AB2").Style = "Normal"
$ws1.Range("AC2").NumberFormat = "@"
$ws1.Range("AC2").Value = "0.0"
$ws1.Range("AC2").Style = "Normal"
$ws1.Range("AD2").NumberFormat = "@"
$ws1.Range("AD2").Value = "1.4"
$ws1.Range("AD2").Style = "Normal"
$ws1.Range("AE2").NumberFormat = "@"
$ws1.Range("AE2").Value = "186"
$ws1.Range("AE2").Style = "Normal"
$ws1.Range("AF2").NumberFormat = "@"
$ws1.Range("AF2").Value = "5.4"
$ws1.Range("AF2").Style = "Normal"
$ws1.Range("AG2").NumberFormat = "@"
$ws1.Range("AG2").Value = "1024.2"
$ws1.Range("AG2").Style = "Normal"
$ws1.Range("AH2").NumberFormat = "@"
$ws1.Range("AH2").Value = "437"
$ws1.Range("AH2").Style = "Normal"
$ws1.Range("AI2").NumberFormat = "@"
$ws1.Range("AI2").Value = "09:30 - 10:00"
$ws1.Range("AI2").Style = "Normal"
$ws1.Range("AJ2").NumberFormat = "@"
$ws1.Range("AJ2").Value = "12.8"
$ws1.Range("AJ2").Style = "Normal"
$ws1.Range("AK2").NumberFormat = "@"
$ws1.Range("AK2").Value = "13.5"
$ws1.Range("AK2").Style = "Normal"
$ws1.Range("AL2").NumberFormat = "@"
$ws1.Range("AL2").Value = "11.9"
$ws1.Range("AL2").Style = "Normal"
$ws1.Range("AM2").NumberFormat = "@"
$ws1.Range("AM2").Value = "50"
$ws1.Range("AM2").Style = "Normal"
$ws1.Range("AN2").NumberFormat = "@"
$ws1.Range("AN2").Value = "0.0"
$ws1.Range("AN2").Style = "Normal"
$ws1.Range("AO2").NumberFormat = "@"
$ws1.Range("AO2").Value = "1.4"
$ws1.Range("AO2").Style = "Normal"
$ws1.Range("AP2").NumberFormat = "@"
$ws1.Range("AP2").Value = "186"
$ws1.Range("AP2").Style = "Normal"
$ws1.Range("AQ2").NumberFormat = "@"
$ws1.Range("AQ2").Value = "5.4"
$ws1.Range("AQ2").Style = "Normal"
$ws1.Range("AR2").NumberFormat = "@"
$ws1.Range("AR2").Value = "1024.2"
$ws1.Range("AR2").Style = "Normal"
$ws1.Range("AS2").NumberFormat = "@"
$ws1.Range("AS2").Value = "437"
$ws1.Range("AS2").Style = "Normal"
$ws1.Range("H3").NumberFormat = "@"
$ws1.Range("H3").Value = "2026-02-20 10:10:57"
$ws1.Range("H3").Style = "Normal"
$ws1.Range("H4").NumberFormat = "@"
$ws1.Range("H4").Value = "2026-02-20 10:10:57"
$ws1.Range("H4").Style = "Normal"
$ws1.Range("H5").NumberFormat = "@"
$ws1.Range("H5").Value = "2026-02-20 10:10:57"
$ws1.Range("H5").Style = "Normal"
$ws1.Range("H6").NumberFormat = "@"
$ws1.Range("H6").Value = "2026-02-20 10:10:57"
$ws1.Range("H6").Style = "Normal"

# --- Sheet 2: Estudi_Capçaleres ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("F2").NumberFormat = "@"
$ws2.Range("F2").Value = "https://www.meteo.cat/observacions/xema/dades?codi=XJ&dia=2026-02-20T09:30Z"
$ws2.Range("F2").Style = "Normal"

Write-Host "Edit complete"
